$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$body = $s.Shapes.Item(4)
$tr = $body.TextFrame.TextRange

# ---------------------------------------------------------------------------
# 1) "Day 3 - Binding Track Activity (show / hide 'Loading' )"
#      -> "Day 3 - Binding Track Activity (show / hide 'Loading' ), Scan Operator"
#    (split the run so the trailing "' )" becomes its own run, then append the
#    new ", Scan Operator" text after the closing paren)
# ---------------------------------------------------------------------------
$day3 = $tr.Paragraphs(4, 1)
$day3Start = $day3.Start
$day3Text = $day3.Text

$afterLoadingRelIdx = $day3Text.IndexOf("Loading") + 7   # right after "...Loading", before the closing quote
$quoteSpace = $tr.Characters($day3Start + $afterLoadingRelIdx, 2)   # the quote + following space
$quoteSpace.Text = $quoteSpace.Text                    # force a clean run split, no text change

$closeParenRelIdx = $day3Text.IndexOf(")")
$closeParen = $tr.Characters($day3Start + $closeParenRelIdx, 1)
$closeParen.InsertAfter(", Scan Operator")

# ---------------------------------------------------------------------------
# 2) "Day 4 - Adding a Reactive Extension to Custom UI Element, "
#    Merge the standalone leading-space run with the following text run
#    (no visible text change, just a run-structure normalization)
# ---------------------------------------------------------------------------
$day4 = $tr.Paragraphs(5, 1)
$day4Start = $day4.Start
$day4Text = $day4.Text

$addingRelIdx = $day4Text.IndexOf(" Adding")
$addingSub = $tr.Characters($day4Start + $addingRelIdx, ($day4Text.Length - $addingRelIdx - 1))
$addingSub.Text = $addingSub.Text

# ---------------------------------------------------------------------------
# 3) "Day " + "5 " -> single run "Day 5 " (still followed by " - Schedulers (...")
# ---------------------------------------------------------------------------
$day5 = $tr.Paragraphs(7, 1)
$day5Start = $day5.Start

$dayWord = $tr.Characters($day5Start, 4)                # "Day "
$dayWord.Delete()

$fiveRun = $tr.Characters($day5Start, 2)                # formerly "5 ", now right after deletion
$fiveRun.Text = "Day 5 "
